$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header row + all data rows)
$lastRow = $ws.UsedRange.Rows.Count

# Swap the contents of columns A (ADM1_EN) and B (ADM1_TH), header included,
# so that A holds the Thai names and B holds the English names.
$colA = $ws.Range("A1:A$lastRow").Value2
$colB = $ws.Range("B1:B$lastRow").Value2

$ws.Range("A1:A$lastRow").Value2 = $colB
$ws.Range("B1:B$lastRow").Value2 = $colA

# Swap the column widths of A and B to match their (now swapped) content
$widthA = $ws.Columns.Item(1).ColumnWidth
$widthB = $ws.Columns.Item(2).ColumnWidth

$ws.Columns.Item(1).ColumnWidth = $widthB
$ws.Columns.Item(2).ColumnWidth = $widthA

# Move the active selection to C3
$ws.Range("C3").Select() | Out-Null
